# Add two new columns, I ("I0") and J ("IF"), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1), matching the style already used by the other header cells (e.g. H1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-39: column I ("I0") and column J ("IF") ---
$data = @(
    @(2, 1, 3),
    @(3, 1, 6),
    @(4, 1, 5),
    @(5, 1, 6),
    @(6, 1, 6),
    @(7, 1, 3),
    @(8, 1, 6),
    @(9, 1, 6),
    @(10, 1, 6),
    @(11, 1, 5),
    @(12, 1, 5),
    @(13, 1, 6),
    @(14, 1, 5),
    @(15, 1, 5),
    @(16, 1, 4),
    @(17, 1, 6),
    @(18, 1, 5),
    @(19, 1, 6),
    @(20, 1, 6),
    @(21, 1, 5),
    @(22, 1, 6),
    @(23, 1, 6),
    @(24, 1, 5),
    @(25, 1, 4),
    @(26, 4, 7),
    @(27, 1, 3),
    @(28, 3, 5),
    @(29, 7, 7),
    @(30, 7, 7),
    @(31, 4, 4),
    @(32, 2, 3),
    @(33, 5, 5),
    @(34, 7, 7),
    @(35, 6, 7),
    @(36, 8, 9),
    @(37, 3, 3),
    @(38, 8, 8),
    @(39, 8, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $if = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}
